# fix: Add Commercial column to import example file
#
# Adds a new "Commercial" column (M) to the "Worksheet" sheet of the
# import example, with sample email addresses for the two sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M1").Value = "Commercial"
$ws.Range("M2").Value = "admin@test.com"
$ws.Range("M3").Value = "user@test.com"
